$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record needs to be inserted as row 178 (pushing the
# existing rows 178-247 down to 179-248). The new record re-uses the
# same "template" values (market, region code, category, quality,
# volume, unit, classification, etc.) as the row it is inserted above,
# only the date, min/avg price, origin and $/Kg columns differ.

$templateRow = 178
$colCount = 18

# Capture the template row's current values before we shift anything.
$rowVals = @()
for ($c = 1; $c -le $colCount; $c++) {
    $rowVals += $ws.Cells.Item($templateRow, $c).Value2
}

# Insert a blank row at 178, shifting 178:247 down to 179:248.
$ws.Rows.Item($templateRow).Insert()

# Re-populate the newly inserted row with the template's values.
for ($c = 1; $c -le $colCount; $c++) {
    $ws.Cells.Item($templateRow, $c).Value2 = $rowVals[$c - 1]
}

# Apply the new record's own values.
$ws.Range("D178").Value2 = 44845
$ws.Range("K178").Value2 = 8500
$ws.Range("M178").Value2 = 8750
$ws.Range("O178").Value2 = "Provincia del Elquí"
$ws.Range("P178").Value2 = 1458
